# Simulated Wild Card round and logged it.
# Appends this game's play-by-play yardage log to the season-long shared
# strings on YDS / ST, and bumps the cumulative season totals on
# OFF / DEF / ST / TURNS / PEN to include the new game.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: per-play yardage logs (space separated numbers), one
# new game's worth of plays appended to each of the four lists.
# ---------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Text + " 8 1 2 3 0 10 -2 6 5 2 8 23 1 13 8"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Text + " 11 33 6 9 11 4 7 10 3 18 4 10 8 21 9 12"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Text + " 3 -1 1 2 1 4 2 4 1 12 8 2 0 26 4 1 0 2 4 5 1 6 0 6 3 -1 6"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Text + " 17 6 12 16 8 15 7 2 9 1 11 6 18 29 7 10 14 12 8 1 6 5 3"

# ---------------------------------------------------------------
# OFF sheet: season offensive totals (Home row 2 / Road row 3)
# ---------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 431
$wsOFF.Range("E2").Value = 20
$wsOFF.Range("F2").Value = 163
$wsOFF.Range("G2").Value = 126
$wsOFF.Range("J2").Value = 78
$wsOFF.Range("N2").Value = 31
$wsOFF.Range("O2").Value = 57
$wsOFF.Range("P2").Value = 35

$wsOFF.Range("B3").Value = 22
$wsOFF.Range("C3").Value = 364
$wsOFF.Range("E3").Value = 71
$wsOFF.Range("F3").Value = 193
$wsOFF.Range("G3").Value = 41
$wsOFF.Range("H3").Value = 54
$wsOFF.Range("I3").Value = 120
$wsOFF.Range("J3").Value = 97
$wsOFF.Range("L3").Value = 552
$wsOFF.Range("M3").Value = 364
$wsOFF.Range("Q3").Value = 1090

# ---------------------------------------------------------------
# DEF sheet: season defensive totals (Home row 2 / Road row 3)
# ---------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 10
$wsDEF.Range("C2").Value = 389
$wsDEF.Range("F2").Value = 120
$wsDEF.Range("G2").Value = 114
$wsDEF.Range("I2").Value = 14
$wsDEF.Range("J2").Value = 64
$wsDEF.Range("N2").Value = 35
$wsDEF.Range("O2").Value = 35
$wsDEF.Range("P2").Value = 18

$wsDEF.Range("B3").Value = 17
$wsDEF.Range("C3").Value = 409
$wsDEF.Range("D3").Value = 9
$wsDEF.Range("E3").Value = 59
$wsDEF.Range("F3").Value = 232
$wsDEF.Range("G3").Value = 78
$wsDEF.Range("H3").Value = 51
$wsDEF.Range("I3").Value = 132
$wsDEF.Range("J3").Value = 113
$wsDEF.Range("L3").Value = 600
$wsDEF.Range("M3").Value = 397
$wsDEF.Range("Q3").Value = 1019

# ---------------------------------------------------------------
# ST sheet: special teams totals + KO/PT distance-return logs
# ---------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 192
$wsST.Range("D2").Value = 108
$wsST.Range("H2").Value = 14
$wsST.Range("I2").Value = 7
$wsST.Range("L2").Value = 36
$wsST.Range("M2").Value = 27

$wsST.Range("B3").Value = 93

$wsST.Range("B6").Value = $wsST.Range("B6").Text + " 40 18 41"
$wsST.Range("D3").Value = $wsST.Range("D3").Text + " 45 33"
$wsST.Range("D4").Value = $wsST.Range("D4").Text + " 0 0"
$wsST.Range("D5").Value = $wsST.Range("D5").Text + " 0 0 7"

# ---------------------------------------------------------------
# TURNS sheet: turnovers (Road row 3)
# ---------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("B3").Value = 6
$wsTURNS.Range("D3").Value = 13
$wsTURNS.Range("E3").Value = 17

# ---------------------------------------------------------------
# PEN sheet: penalties
# ---------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B3").Value = 44
$wsPEN.Range("D4").Value = 18
